# "final for 2E python section"
#
# The source data was regenerated so that:
#   - "Argentina price price" + the blank spacer column ("_1") collapse
#     into a single merged header "Argentina_priceprice"
#   - "Argentina points points" becomes "Argentina_pointspoints" and moves
#     from column D into column C (the spacer column C is removed)
#
# Net effect on the worksheet: the empty spacer column C is deleted, the old
# column D (Argentina points/points data) becomes the new column C, and the
# two header cells get their new merged text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "price" header text (column B keeps its position)
$ws.Range("B1").Value = "Argentina_priceprice"

# Remove the empty spacer column C - this shifts the old column D
# (the "Argentina points points" data) left into column C
$ws.Columns.Item(3).Delete() | Out-Null

# Update the "points" header text, now living in column C
$ws.Range("C1").Value = "Argentina_pointspoints"

# Re-apply column widths for the new layout (closest achievable widths;
# target character widths are 22.7109375 / 24.7109375)
$ws.Columns.Item(2).ColumnWidth = 21.75
$ws.Columns.Item(3).ColumnWidth = 23.75
